$d = $word.ActiveDocument

# 1) Remove the blank paragraph that sits between the "<Assessment of
#    Significance>" paragraph and the "Having considered..." paragraph —
#    i.e. merge the empty paragraph into the following one by deleting
#    the paragraph mark that ends the empty paragraph.
$found = $d.Content.Find.Execute("Having considered the submitted document")
$hitPara = $d.Content.Find.Parent.Paragraphs(1)

# Locate paragraph containing "Having considered"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Having considered the submitted document*") {
        $havingPara = $p
        break
    }
}
$prevPara = $havingPara.Previous()
if ($prevPara.Range.Text.Trim() -eq "") {
    $prevPara.Range.Delete()
}

# 2) Split the "results of the archaeological investigation..." run so the
#    lastRenderedPageBreak marker sits right before "Record." instead of
#    at the very start of the paragraph.
$d.Content.Find.Execute("Environment Record. The site archive", $true, $false, $false, $false, $false, $true, 1, $false, "Environment `rRecord. The site archive", 2)
